$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 12.056684
$ws.Range("H2").Value = 36.170052
$ws.Range("I2").Value = 0.06307822458376462
$ws.Range("J2").Value = 0.06307822458376462
$ws.Range("M2").Value = 3.759736666666667
$ws.Range("N2").Value = 11.27921
$ws.Range("O2").Value = 0.0683751702595819
$ws.Range("P2").Value = 0.06837517025958188
$ws.Range("Q2").Value = 45.32995691321333
$ws.Range("R2").Value = 407.96961221892
$ws.Range("S2").Value = 0.00431298434558705
$ws.Range("T2").Value = 0.004312984345587049

$ws.Range("G3").Value = 12.056684
$ws.Range("H3").Value = 36.170052
$ws.Range("I3").Value = 0.06307822458376462
$ws.Range("J3").Value = 0.06307822458376462
$ws.Range("O3").Value = 0.6514180024294648
$ws.Range("P3").Value = 0.6514180024294647
$ws.Range("Q3").Value = 431.8636410046973
$ws.Range("R3").Value = 3886.772769042276
$ws.Range("S3").Value = 0.04109029105515311
$ws.Range("T3").Value = 0.0410902910551531

$ws.Range("G4").Value = 12.056684
$ws.Range("H4").Value = 36.170052
$ws.Range("I4").Value = 0.06307822458376462
$ws.Range("J4").Value = 0.06307822458376462
$ws.Range("O4").Value = 0.2802068273109533
$ws.Range("P4").Value = 0.2802068273109533
$ws.Range("Q4").Value = 185.7657298778533
$ws.Range("R4").Value = 1671.89156890068
$ws.Range("S4").Value = 0.01767494918302446
$ws.Range("T4").Value = 0.01767494918302446

$ws.Range("I5").Value = 0.1315309049843414
$ws.Range("J5").Value = 0.1315309049843414
$ws.Range("M5").Value = 3.759736666666667
$ws.Range("N5").Value = 11.27921
$ws.Range("O5").Value = 0.0683751702595819
$ws.Range("P5").Value = 0.06837517025958188
$ws.Range("Q5").Value = 94.52216347304667
$ws.Range("R5").Value = 850.69947125742
$ws.Range("S5").Value = 0.00899344802270123
$ws.Range("T5").Value = 0.008993448022701227

$ws.Range("I6").Value = 0.1315309049843414
$ws.Range("J6").Value = 0.1315309049843414
$ws.Range("O6").Value = 0.6514180024294648
$ws.Range("P6").Value = 0.6514180024294647
$ws.Range("S6").Value = 0.08568159938263939
$ws.Range("T6").Value = 0.08568159938263938

$ws.Range("I7").Value = 0.1315309049843414
$ws.Range("J7").Value = 0.1315309049843414
$ws.Range("O7").Value = 0.2802068273109533
$ws.Range("P7").Value = 0.2802068273109533
$ws.Range("S7").Value = 0.03685585757900075
$ws.Range("T7").Value = 0.03685585757900074

$ws.Range("H8").Value = 461.8238680000001
$ws.Range("I8").Value = 0.8053908704318941
$ws.Range("J8").Value = 0.8053908704318941
$ws.Range("M8").Value = 3.759736666666667
$ws.Range("N8").Value = 11.27921
$ws.Range("O8").Value = 0.0683751702595819
$ws.Range("P8").Value = 0.06837517025958188
$ws.Range("Q8").Value = 578.7787100204756
$ws.Range("R8").Value = 5209.00839018428
$ws.Range("S8").Value = 0.05506873789129362
$ws.Range("T8").Value = 0.05506873789129361

$ws.Range("H9").Value = 461.8238680000001
$ws.Range("I9").Value = 0.8053908704318941
$ws.Range("J9").Value = 0.8053908704318941
$ws.Range("O9").Value = 0.6514180024294648
$ws.Range("P9").Value = 0.6514180024294647
$ws.Range("S9").Value = 0.5246461119916723
$ws.Range("T9").Value = 0.5246461119916722

$ws.Range("H10").Value = 461.8238680000001
$ws.Range("I10").Value = 0.8053908704318941
$ws.Range("J10").Value = 0.8053908704318941
$ws.Range("O10").Value = 0.2802068273109533
$ws.Range("P10").Value = 0.2802068273109533
$ws.Range("S10").Value = 0.2256760205489281
$ws.Range("T10").Value = 0.2256760205489281
